# The name "Trần Thanh Vy" (stored as two runs: "Trần" and " Thanh Vy",
# the second one tagged w:lang="vi-VN") shows up in two different rows of
# the single big table in this journal document. The diff only touches
# the one belonging to the "3.Chuẩn bị các file báo cáo đồ án." task row
# (week 02) -- the other "Trần Thanh Vy" cell (week 01 task "2.4. Lên kế
# hoạch thực hiện ban đầu") must stay untouched.
#
# We locate that specific cell by matching the task-description text in
# the neighbouring column, then merge the two runs into a single run
# whose text is "Mai Nhật Hào" (keeping the first run's formatting, which
# has no vi-VN language tag -- matching the target XML).

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$targetTask = "3.Chuẩn bị các file báo cáo đồ án."
$oldFirstRunText = "Trần"
$newText = "Mai Nhật Hào"

for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $row = $table.Rows.Item($r)
    if ($row.Cells.Count -lt 3) { continue }

    $taskCell = $row.Cells.Item(2)
    $taskText = $taskCell.Range.Text.TrimEnd([char]13, [char]7)

    if ($taskText -eq $targetTask) {
        $nameCell = $row.Cells.Item(3)
        $nameRange = $nameCell.Range
        $fullText = $nameRange.Text.TrimEnd([char]13, [char]7)

        if ($fullText -eq "Trần Thanh Vy") {
            $cellStart = $nameRange.Start
            $boundary = $cellStart + $oldFirstRunText.Length
            $cellTextEnd = $cellStart + $fullText.Length

            # Delete the trailing " Thanh Vy" run (2nd run) first so the
            # first run's range stays valid.
            $tailRange = $d.Range($boundary, $cellTextEnd)
            $tailRange.Delete()

            # Replace the remaining "Trần" run's text with the new name.
            $headRange = $d.Range($cellStart, $boundary)
            $headRange.Text = $newText
        }
        break
    }
}
